$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(45370, 10, 4.220987319946289),
    @(45370, 11, 3.964477062225342),
    @(45370, 12, 4.176630973815918),
    @(45370, 13, 4.159878730773926),
    @(45370, 14, 3.80241584777832),
    @(45370, 15, 3.707133531570435),
    @(45370, 16, 3.794259309768677),
    @(45370, 17, 3.639703035354614),
    @(45370, 18, 3.640395879745483),
    @(45370, 19, 3.430887699127197),
    @(45370, 20, 3.341203212738037),
    @(45370, 21, 3.365144968032837),
    @(45370, 22, 2.86371636390686),
    @(45370, 23, 2.356311798095703),
    @(45371, 0, 2.274955749511719),
    @(45371, 1, 2.265320539474487),
    @(45371, 2, 2.13236403465271),
    @(45371, 3, 2.090662717819214),
    @(45371, 4, 2.152424573898315),
    @(45371, 5, 2.757708787918091),
    @(45371, 6, 3.065781831741333),
    @(45371, 7, 4.204745769500732),
    @(45371, 8, 4.178340911865234),
    @(45371, 9, 4.327641010284424),
    @(45371, 10, 4.12451171875),
    @(45371, 11, 3.838339328765869),
    @(45371, 12, 4.091320514678955),
    @(45371, 13, 3.999934434890747),
    @(45371, 14, 3.738850593566895),
    @(45371, 15, 3.457412004470825),
    @(45371, 16, 3.672144174575806),
    @(45371, 17, 3.706796646118164),
    @(45371, 18, 3.697768926620483),
    @(45371, 19, 3.300342082977295),
    @(45371, 20, 3.369502305984497),
    @(45371, 21, 3.574325323104858),
    @(45371, 22, 2.986258029937744),
    @(45371, 23, 2.339095830917358),
    @(45372, 0, 2.29724383354187),
    @(45372, 1, 2.287608623504639),
    @(45372, 2, 2.170822858810425),
    @(45372, 3, 2.105140686035156),
    @(45372, 4, 2.165964126586914),
    @(45372, 5, 2.736754655838013),
    @(45372, 6, 3.088069915771484),
    @(45372, 7, 4.221695423126221),
    @(45372, 8, 4.296972274780273),
    @(45372, 9, 4.318789005279541),
    @(45372, 10, 4.155550956726074),
    @(45372, 11, 3.979258298873901),
    @(45372, 12, 4.177427768707275),
    @(45372, 13, 4.174818992614746),
    @(45372, 14, 3.795731067657471),
    @(45372, 15, 3.769979476928711),
    @(45372, 16, 3.720515251159668),
    @(45372, 17, 3.68045973777771),
    @(45372, 18, 3.625741958618164),
    @(45372, 19, 3.428699016571045),
    @(45372, 20, 3.479614973068237),
    @(45372, 21, 3.640857934951782),
    @(45372, 22, 3.012799978256226),
    @(45372, 23, 2.36231541633606),
    @(45373, 0, 2.24016261100769),
    @(45373, 1, 2.209706544876099),
    @(45373, 2, 2.113217353820801),
    @(45373, 3, 2.070497035980225),
    @(45373, 4, 2.128940582275391),
    @(45373, 5, 2.734561204910278),
    @(45373, 6, 3.064020156860352),
    @(45373, 7, 4.280899524688721),
    @(45373, 8, 4.297880172729492),
    @(45373, 9, 4.29173755645752),
    @(45373, 10, 4.162452220916748),
    @(45373, 11, 3.807162761688232),
    @(45373, 12, 3.971662521362305),
    @(45373, 13, 3.969053745269775),
    @(45373, 14, 3.71341347694397),
    @(45373, 15, 3.764338731765747),
    @(45373, 16, 3.776564598083496),
    @(45373, 17, 3.710996866226196),
    @(45373, 18, 3.613184452056885),
    @(45373, 19, 3.620223522186279),
    @(45373, 20, 3.545651435852051),
    @(45373, 21, 3.424517631530762),
    @(45373, 22, 2.705356359481812),
    @(45373, 23, 2.370362281799316),
    @(45374, 0, 2.15131402015686),
    @(45374, 1, 2.246299028396606),
    @(45374, 2, 2.149628877639771),
    @(45374, 3, 2.109325408935547),
    @(45374, 4, 1.993576765060425),
    @(45374, 5, 2.128603219985962),
    @(45374, 6, 1.990267634391785),
    @(45374, 7, 1.805531144142151),
    @(45374, 8, 1.789688348770142),
    @(45374, 9, 1.463794946670532)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $i + 2
    $row = $values[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}